# Regenerate orders with updated distance/size codes.
# Mapping applied to any textual cell content:
#   D64 -> D69
#   D51 -> D55
#   D80 -> D86
#   S30 -> S31
# (S20 and S25 are left unchanged.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ur = $ws.UsedRange
$firstRow = $ur.Row
$firstCol = $ur.Column
$lastRow = $firstRow + $ur.Rows.Count - 1
$lastCol = $firstCol + $ur.Columns.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($val -ne $null -and $val.GetType().FullName -eq "System.String") {
            $newVal = $val
            $newVal = $newVal.Replace("D64", "D69")
            $newVal = $newVal.Replace("D51", "D55")
            $newVal = $newVal.Replace("D80", "D86")
            $newVal = $newVal.Replace("S30", "S31")
            if ($newVal -ne $val) {
                $cell.Value2 = $newVal
            }
        }
    }
}
